$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00609756097560976
$ws.Range("C2").Value = 0.981707317073171
$ws.Range("D2").Value = 0.00203252032520325
$ws.Range("E2").Value = 0.00203252032520325
$ws.Range("F2").Value = 0.016260162601626
$ws.Range("G2").Value = 0.00203252032520325
$ws.Range("H2").Value = 0.0121951219512195
$ws.Range("I2").Value = 0.973577235772358
$ws.Range("J2").Value = 0.024390243902439
$ws.Range("K2").Value = 0.951219512195122
$ws.Range("L2").Value = 0.961382113821138
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0.0040650406504065
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.995934959349594
$ws.Range("R2").Value = 0.0040650406504065
$ws.Range("S2").Value = 0.890243902439024
$ws.Range("T2").Value = 0.0609756097560976
$ws.Range("U2").Value = 0.0691056910569106
$ws.Range("V2").Value = 0.989837398373984
$ws.Range("W2").Value = 0.00609756097560976
$ws.Range("X2").Value = 0.00203252032520325

$ws.Range("B3").Value = 0.0101626016260163
$ws.Range("C3").Value = 0.0040650406504065
$ws.Range("D3").Value = 0.0223577235772358
$ws.Range("E3").Value = 0.00609756097560976
$ws.Range("F3").Value = 0.00203252032520325
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.977642276422764
$ws.Range("I3").Value = 0.0101626016260163
$ws.Range("J3").Value = 0.0142276422764228
$ws.Range("K3").Value = 0.00203252032520325
$ws.Range("L3").Value = 0.0040650406504065
$ws.Range("M3").Value = 0.943089430894309
$ws.Range("N3").Value = 0.00609756097560976
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.991869918699187
$ws.Range("Q3").Value = 0.0040650406504065
$ws.Range("R3").Value = 0.99390243902439
$ws.Range("S3").Value = 0.0975609756097561
$ws.Range("T3").Value = 0.930894308943089
$ws.Range("U3").Value = 0.92479674796748
$ws.Range("V3").Value = 0.00203252032520325
$ws.Range("W3").Value = 0.0040650406504065
$ws.Range("X3").Value = 0.0040650406504065

$ws.Range("B4").Value = 0.979674796747967
$ws.Range("C4").Value = 0.00203252032520325
$ws.Range("D4").Value = 0.00813008130081301
$ws.Range("E4").Value = 0.0040650406504065
$ws.Range("F4").Value = 0.975609756097561
$ws.Range("G4").Value = 0.99390243902439
$ws.Range("H4").Value = 0.00813008130081301
$ws.Range("I4").Value = 0.0040650406504065
$ws.Range("J4").Value = 0.951219512195122
$ws.Range("K4").Value = 0.0040650406504065
$ws.Range("L4").Value = 0.032520325203252
$ws.Range("M4").Value = 0.00203252032520325
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.00203252032520325
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.00203252032520325
$ws.Range("S4").Value = 0.0101626016260163
$ws.Range("T4").Value = 0.00203252032520325
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0.00609756097560976
$ws.Range("W4").Value = 0.98780487804878
$ws.Range("X4").Value = 0.991869918699187

$ws.Range("B5").Value = 0.0040650406504065
$ws.Range("C5").Value = 0.0121951219512195
$ws.Range("D5").Value = 0.967479674796748
$ws.Range("E5").Value = 0.98780487804878
$ws.Range("F5").Value = 0.00609756097560976
$ws.Range("G5").Value = 0.0040650406504065
$ws.Range("H5").Value = 0.00203252032520325
$ws.Range("I5").Value = 0.0121951219512195
$ws.Range("J5").Value = 0.0101626016260163
$ws.Range("K5").Value = 0.040650406504065
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0.0528455284552846
$ws.Range("N5").Value = 0.989837398373984
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0.00609756097560976
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.00203252032520325
$ws.Range("T5").Value = 0.0040650406504065
$ws.Range("U5").Value = 0.00609756097560976
$ws.Range("V5").Value = 0.00203252032520325
$ws.Range("W5").Value = 0.00203252032520325
$ws.Range("X5").Value = 0.00203252032520325

